$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Z (26th column), shifting existing
# columns Z.. right by one.
$ws.Columns.Item(26).Insert()

# New column header + data for the inserted "STAT" column
$ws.Cells.Item(1, 26).Value = "STAT"
$ws.Cells.Item(2, 26).Value = "I"
$ws.Cells.Item(3, 26).Value = "I"

# Copy style from neighboring column (AA, now shifted from old Z) for header/data rows
$ws.Cells.Item(1, 26).Style = $ws.Cells.Item(1, 27).Style
$ws.Cells.Item(2, 26).Style = $ws.Cells.Item(2, 27).Style
$ws.Cells.Item(3, 26).Style = $ws.Cells.Item(3, 27).Style
$ws.Cells.Item(4, 26).Style = $ws.Cells.Item(4, 27).Style

# Set the narrower column width for the new STAT column
$ws.Columns.Item(26).ColumnWidth = 6.33203125

# Update the view: scroll so column Y is the top-left visible column,
# and set the active selection to AF10
$ws.Application.ActiveWindow.ScrollColumn = 25
$ws.Range("AF10").Select()
